$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("东方日升", "隆基绿能", "锋龙股份"),
    @("钧达股份", "东方日升", "白银有色"),
    @("乾照光电", "锋龙股份", "利欧股份"),
    @("航天电子", "航天电子", "航天电子"),
    @("拓日新能", "白银有色", "隆基绿能"),
    @("锋龙股份", "乾照光电", "金风科技"),
    @("金风科技", "钧达股份", "东方日升"),
    @("隆基绿能", "拓日新能", "湖南白银"),
    @("中国卫星", "利欧股份", "钧达股份"),
    @("白银有色", "金风科技", "乾照光电"),
    @("迈为股份", "协鑫集成", "拓日新能"),
    @("利欧股份", "捷佳伟创", "蓝色光标"),
    @("巨力索具", "湖南白银", "中国卫星"),
    @("湖南白银", "迈为股份", "岩山科技"),
    @("捷佳伟创", "贵州茅台", "协鑫集成"),
    @("协鑫集成", "中国卫星", "平潭发展"),
    @("蓝色光标", "东方财富", "特变电工"),
    @("岩山科技", "铜陵有色", "航天发展"),
    @("铜陵有色", "蓝色光标", "巨力索具"),
    @("臻镭科技", "岩山科技", "铜陵有色")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
